$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.252.31"
$ws.Range("E2").Value = "  -6.03%  "
$ws.Range("D3").Value = "1.670.58"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.94"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5081"
$ws.Range("E6").Value = "  -12.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2661"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06346"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.59"
$ws.Range("E10").Value = "  -6.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07372"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "1.679.60"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.544"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5798"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("D15").Value = "1.897.16"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008543"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.05"
$ws.Range("E17").Value = "  -12.87%  "
$ws.Range("D18").Value = "26.306.87"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.943"
$ws.Range("E19").Value = "  -7.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.84"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.30"
$ws.Range("E22").Value = "  -7.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.198"
$ws.Range("E23").Value = "  -6.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.74"
$ws.Range("E25").Value = "  -4.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.692"
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1172"
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.72"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05802"
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("E30").Value = "  -7.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.323"
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.532"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.516"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.639"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.012"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5983"
$ws.Range("E36").Value = "  -6.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.362"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.642"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.037"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "1.077.19"
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8598"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.83"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "1.823.80"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.78"
$ws.Range("E47").Value = "  -5.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.081"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4300"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05188"
$ws.Range("E51").Value = "  -3.54%  "
